$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so that numeric-looking strings
# (e.g. "0.140", "51.447.45", "6.80") are written verbatim instead of
# being auto-coerced into numbers (which would drop trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.447.45"
$ws.Range("E2").Value = "  +5.02%  "
$ws.Range("D3").Value = "2.730.82"
$ws.Range("E3").Value = "  +4.12%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "115.29"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("D6").Value = "330.98"
$ws.Range("E6").Value = "  +2.63%  "
$ws.Range("D7").Value = "0.536"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("D9").Value = "0.566"
$ws.Range("E9").Value = "  +4.50%  "
$ws.Range("D10").Value = "41.27"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("D12").Value = "20.01"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("D14").Value = "7.55"
$ws.Range("E14").Value = "  +4.35%  "
$ws.Range("D15").Value = "3.169.70"
$ws.Range("E15").Value = "  +4.52%  "
$ws.Range("D16").Value = "2.767.14"
$ws.Range("E16").Value = "  +5.52%  "
$ws.Range("D17").Value = "0.874"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "51.446.09"
$ws.Range("E18").Value = "  +5.04%  "
$ws.Range("D19").Value = "3.12"
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "13.31"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "6.80"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").Value = "278.22"
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("D24").Value = "69.00"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "2.63"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").Value = "26.57"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "0.140"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "49.94"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").Value = "5.51"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "0.0815"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "18.97"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "4.98"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "3.18"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").Value = "127.65"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "2.29"
$ws.Range("E41").Value = "  +7.52%  "
$ws.Range("D42").Value = "0.0343"
$ws.Range("E42").Value = "  +8.21%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "22.89"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("E45").Value = "  +10.93%  "
$ws.Range("D46").Value = "2.082.35"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "2.22"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("E49").Value = "  +6.59%  "
$ws.Range("D50").Value = "8.89"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "59.51"
$ws.Range("E51").Value = "  +1.77%  "

# Restore the original (default/General) number format and style for column D
# now that the text values are safely stored.
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
